$wb = $excel.ActiveWorkbook

# ALC row 41
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 8547156
$ws.Range("I41").Value = 12345887
$ws.Range("J41").Value = 13.5
$ws.Range("K41").Value = 12345887
$ws.Range("L41").Value = 13.5
$ws.Range("M41").Value = -12345447
$ws.Range("N41").Value = -893.5

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 15152852
$ws.Range("I112").Value = 950
$ws.Range("J112").Value = 18183232
$ws.Range("K112").Value = 2850
$ws.Range("L112").Value = 54549696
$ws.Range("M112").Value = -1742
$ws.Range("N112").Value = -54551912

# ALC row 127
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 706.1905
$ws.Range("I127").Value = 383
$ws.Range("J127").Value = 1000
$ws.Range("K127").Value = 1149
$ws.Range("L127").Value = 3000
$ws.Range("M127").Value = 3811
$ws.Range("N127").Value = -12920

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 977.9375
$ws.Range("J129").Value = 1111.3846
$ws.Range("L129").Value = 3334.1538
$ws.Range("N129").Value = -13334.1538

# ALC row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 3333.7058
$ws.Range("I131").Value = 1981.9166
$ws.Range("J131").Value = 6578
$ws.Range("K131").Value = 5945.7498
$ws.Range("L131").Value = 19734
$ws.Range("M131").Value = -905.7497999999996
$ws.Range("N131").Value = -29814

# ALC row 133
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 42693.57
$ws.Range("J133").Value = 42693.57
$ws.Range("L133").Value = 42693.57
$ws.Range("N133").Value = -52813.57

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 45455750
$ws.Range("I137").Value = 66667796
$ws.Range("J137").Value = 1369.1428
$ws.Range("K137").Value = 200003388
$ws.Range("L137").Value = 4107.428400000001
$ws.Range("M137").Value = -200000838
$ws.Range("N137").Value = -9207.428400000001

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 6351808.5
$ws.Range("I138").Value = 2270245.2
$ws.Range("J138").Value = 8549574
$ws.Range("K138").Value = 6810735.600000001
$ws.Range("L138").Value = 25648722
$ws.Range("M138").Value = -6805595.600000001
$ws.Range("N138").Value = -25659002

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21124.807
$ws.Range("I32").Value = 2766.6667
$ws.Range("J32").Value = 177169
$ws.Range("K32").Value = 2766.6667
$ws.Range("L32").Value = 177169
$ws.Range("M32").Value = -2479.6667
$ws.Range("N32").Value = -177743

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6081.92
$ws.Range("I74").Value = 1002.2857
$ws.Range("K74").Value = 1002.2857
$ws.Range("M74").Value = -128.2857

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 6081.92
$ws.Range("I77").Value = 1002.2857
$ws.Range("K77").Value = 5011.4285
$ws.Range("M77").Value = -643.4285

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1208.0741
$ws.Range("I20").Value = 1141.1111
$ws.Range("J20").Value = 1342
$ws.Range("K20").Value = 1141.1111
$ws.Range("L20").Value = 1342
$ws.Range("M20").Value = -894.1111000000001
$ws.Range("N20").Value = -1836

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1619.762
$ws.Range("I107").Value = 1530.5625
$ws.Range("J107").Value = 1905.2
$ws.Range("K107").Value = 1530.5625
$ws.Range("L107").Value = 1905.2
$ws.Range("M107").Value = 389.4375
$ws.Range("N107").Value = -5745.2

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4312.925
$ws.Range("I31").Value = 1107.0344
$ws.Range("J31").Value = 12764.818
$ws.Range("K31").Value = 1107.0344
$ws.Range("L31").Value = 12764.818
$ws.Range("M31").Value = -812.0344
$ws.Range("N31").Value = -13354.818

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4312.925
$ws.Range("I34").Value = 1107.0344
$ws.Range("J34").Value = 12764.818
$ws.Range("K34").Value = 1107.0344
$ws.Range("L34").Value = 12764.818
$ws.Range("M34").Value = -905.0344
$ws.Range("N34").Value = -13168.818

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 400.2381
$ws.Range("I107").Value = 230.55556
$ws.Range("J107").Value = 527.5
$ws.Range("K107").Value = 230.55556
$ws.Range("L107").Value = 527.5
$ws.Range("M107").Value = 1689.44444
$ws.Range("N107").Value = -4367.5

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1819.9286
$ws.Range("I122").Value = 1568.4286
$ws.Range("J122").Value = 2071.4285
$ws.Range("K122").Value = 4705.2858
$ws.Range("L122").Value = 6214.2855
$ws.Range("M122").Value = -2255.2858
$ws.Range("N122").Value = -11114.2855

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2623.4707
$ws.Range("I134").Value = 1398.4
$ws.Range("J134").Value = 6026.4443
$ws.Range("K134").Value = 4195.200000000001
$ws.Range("L134").Value = 18079.3329
$ws.Range("M134").Value = -1660.200000000001
$ws.Range("N134").Value = -23149.3329

# CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 512.1667
$ws.Range("I107").Value = 521.3333
$ws.Range("J107").Value = 484.66666
$ws.Range("K107").Value = 1563.9999
$ws.Range("L107").Value = 1453.99998
$ws.Range("M107").Value = 356.0001
$ws.Range("N107").Value = -5293.999980000001

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 924.8
$ws.Range("I122").Value = 262.5
$ws.Range("J122").Value = 1366.3334
$ws.Range("K122").Value = 2362.5
$ws.Range("L122").Value = 12297.0006
$ws.Range("M122").Value = 87.5
$ws.Range("N122").Value = -17197.0006

# GSM row 40
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 8000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 8000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 8000
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -8302

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2788.889
$ws.Range("J126").Value = 2868.4211
$ws.Range("L126").Value = 8605.263300000001
$ws.Range("N126").Value = -13545.2633

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3582.64
$ws.Range("I132").Value = 3047.2354
$ws.Range("J132").Value = 4720.375
$ws.Range("K132").Value = 9141.706200000001
$ws.Range("L132").Value = 14161.125
$ws.Range("M132").Value = -6611.706200000001
$ws.Range("N132").Value = -19221.125

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2100.4
$ws.Range("I68").Value = 1714.8572
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 1714.8572
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -965.8571999999999
$ws.Range("N68").Value = -4498

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2100.4
$ws.Range("I71").Value = 1714.8572
$ws.Range("J71").Value = 3000
$ws.Range("K71").Value = 8574.286
$ws.Range("L71").Value = 15000
$ws.Range("M71").Value = -4830.286
$ws.Range("N71").Value = -22488

# WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 428.18182
$ws.Range("I113").Value = 336
$ws.Range("J113").Value = 505
$ws.Range("K113").Value = 1008
$ws.Range("L113").Value = 1515
$ws.Range("M113").Value = 1162
$ws.Range("N113").Value = -5855

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 40682.16
$ws.Range("I126").Value = 53181.26
$ws.Range("J126").Value = 1101.6666
$ws.Range("K126").Value = 159543.78
$ws.Range("L126").Value = 3304.9998
$ws.Range("M126").Value = -157073.78
$ws.Range("N126").Value = -8244.9998

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2531.5134
$ws.Range("J132").Value = 2864.125
$ws.Range("L132").Value = 8592.375
$ws.Range("N132").Value = -13652.375
